$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 110.8604276666667
$ws.Cells.Item(2, 8).Value = 332.581283
$ws.Cells.Item(2, 9).Value = 0.2509786052589675
$ws.Cells.Item(2, 10).Value = 0.2509786052589675
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 23.63579766666667
$ws.Cells.Item(2, 14).Value = 70.907393
$ws.Cells.Item(2, 15).Value = 0.06827844587621175
$ws.Cells.Item(2, 16).Value = 0.06827844587621175
$ws.Cells.Item(2, 17).Value = 2620.274637569469
$ws.Cells.Item(2, 18).Value = 23582.47173812522
$ws.Cells.Item(2, 19).Value = 0.01713642911526152
$ws.Cells.Item(2, 20).Value = 0.01713642911526152

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 110.8604276666667
$ws.Cells.Item(3, 8).Value = 332.581283
$ws.Cells.Item(3, 9).Value = 0.2509786052589675
$ws.Cells.Item(3, 10).Value = 0.2509786052589675
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 181.2883913333334
$ws.Cells.Item(3, 14).Value = 543.865174
$ws.Cells.Item(3, 15).Value = 0.5237009467675041
$ws.Cells.Item(3, 16).Value = 0.523700946767504
$ws.Cells.Item(3, 17).Value = 20097.70859421536
$ws.Cells.Item(3, 18).Value = 180879.3773479383
$ws.Cells.Item(3, 19).Value = 0.131437733192509
$ws.Cells.Item(3, 20).Value = 0.1314377331925089

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 110.8604276666667
$ws.Cells.Item(4, 8).Value = 332.581283
$ws.Cells.Item(4, 9).Value = 0.2509786052589675
$ws.Cells.Item(4, 10).Value = 0.2509786052589675
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 111.1005463333333
$ws.Cells.Item(4, 14).Value = 333.301639
$ws.Cells.Item(4, 15).Value = 0.3209442197221123
$ws.Cells.Item(4, 16).Value = 0.3209442197221123
$ws.Cells.Item(4, 17).Value = 12316.65408051365
$ws.Cells.Item(4, 18).Value = 110849.8867246228
$ws.Cells.Item(4, 19).Value = 0.08055013263178334
$ws.Cells.Item(4, 20).Value = 0.08055013263178334

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 110.8604276666667
$ws.Cells.Item(5, 8).Value = 332.581283
$ws.Cells.Item(5, 9).Value = 0.2509786052589675
$ws.Cells.Item(5, 10).Value = 0.2509786052589675
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 30.14303933333333
$ws.Cells.Item(5, 14).Value = 90.42911799999999
$ws.Cells.Item(5, 15).Value = 0.08707638763417187
$ws.Cells.Item(5, 16).Value = 0.08707638763417187
$ws.Cells.Item(5, 17).Value = 3341.670231666488
$ws.Cells.Item(5, 18).Value = 30075.03208499839
$ws.Cells.Item(5, 19).Value = 0.02185431031941366
$ws.Cells.Item(5, 20).Value = 0.02185431031941366

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 184.841802
$ws.Cells.Item(6, 8).Value = 554.525406
$ws.Cells.Item(6, 9).Value = 0.4184661617850055
$ws.Cells.Item(6, 10).Value = 0.4184661617850055
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 23.63579766666667
$ws.Cells.Item(6, 14).Value = 70.907393
$ws.Cells.Item(6, 15).Value = 0.06827844587621175
$ws.Cells.Item(6, 16).Value = 0.06827844587621175
$ws.Cells.Item(6, 17).Value = 4368.883432414062
$ws.Cells.Item(6, 18).Value = 39319.95089172656
$ws.Cells.Item(6, 19).Value = 0.02857221917846356
$ws.Cells.Item(6, 20).Value = 0.02857221917846356

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 184.841802
$ws.Cells.Item(7, 8).Value = 554.525406
$ws.Cells.Item(7, 9).Value = 0.4184661617850055
$ws.Cells.Item(7, 10).Value = 0.4184661617850055
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 181.2883913333334
$ws.Cells.Item(7, 14).Value = 543.865174
$ws.Cells.Item(7, 15).Value = 0.5237009467675041
$ws.Cells.Item(7, 16).Value = 0.523700946767504
$ws.Cells.Item(7, 17).Value = 33509.67293573452
$ws.Cells.Item(7, 18).Value = 301587.0564216106
$ws.Cells.Item(7, 19).Value = 0.2191511251169709
$ws.Cells.Item(7, 20).Value = 0.2191511251169709

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 184.841802
$ws.Cells.Item(8, 8).Value = 554.525406
$ws.Cells.Item(8, 9).Value = 0.4184661617850055
$ws.Cells.Item(8, 10).Value = 0.4184661617850055
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 111.1005463333333
$ws.Cells.Item(8, 14).Value = 333.301639
$ws.Cells.Item(8, 15).Value = 0.3209442197221123
$ws.Cells.Item(8, 16).Value = 0.3209442197221123
$ws.Cells.Item(8, 17).Value = 20536.02518743783
$ws.Cells.Item(8, 18).Value = 184824.2266869404
$ws.Cells.Item(8, 19).Value = 0.1343042957741958
$ws.Cells.Item(8, 20).Value = 0.1343042957741958

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 184.841802
$ws.Cells.Item(9, 8).Value = 554.525406
$ws.Cells.Item(9, 9).Value = 0.4184661617850055
$ws.Cells.Item(9, 10).Value = 0.4184661617850055
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 30.14303933333333
$ws.Cells.Item(9, 14).Value = 90.42911799999999
$ws.Cells.Item(9, 15).Value = 0.08707638763417187
$ws.Cells.Item(9, 16).Value = 0.08707638763417187
$ws.Cells.Item(9, 17).Value = 5571.693708130211
$ws.Cells.Item(9, 18).Value = 50145.2433731719
$ws.Cells.Item(9, 19).Value = 0.03643852171537521
$ws.Cells.Item(9, 20).Value = 0.03643852171537521

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 95.23175666666667
$ws.Cells.Item(10, 8).Value = 285.69527
$ws.Cells.Item(10, 9).Value = 0.2155966197102082
$ws.Cells.Item(10, 10).Value = 0.2155966197102082
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 23.63579766666667
$ws.Cells.Item(10, 14).Value = 70.907393
$ws.Cells.Item(10, 15).Value = 0.06827844587621175
$ws.Cells.Item(10, 16).Value = 0.06827844587621175
$ws.Cells.Item(10, 17).Value = 2250.878532014568
$ws.Cells.Item(10, 18).Value = 20257.90678813111
$ws.Cells.Item(10, 19).Value = 0.01472060212997766
$ws.Cells.Item(10, 20).Value = 0.01472060212997766

# Row 11
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 95.23175666666667
$ws.Cells.Item(11, 8).Value = 285.69527
$ws.Cells.Item(11, 9).Value = 0.2155966197102082
$ws.Cells.Item(11, 10).Value = 0.2155966197102082
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 181.2883913333334
$ws.Cells.Item(11, 14).Value = 543.865174
$ws.Cells.Item(11, 15).Value = 0.5237009467675041
$ws.Cells.Item(11, 16).Value = 0.523700946767504
$ws.Cells.Item(11, 17).Value = 17264.41196994744
$ws.Cells.Item(11, 18).Value = 155379.707729527
$ws.Cells.Item(11, 19).Value = 0.1129081538621096
$ws.Cells.Item(11, 20).Value = 0.1129081538621096

# Row 12
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 95.23175666666667
$ws.Cells.Item(12, 8).Value = 285.69527
$ws.Cells.Item(12, 9).Value = 0.2155966197102082
$ws.Cells.Item(12, 10).Value = 0.2155966197102082
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 111.1005463333333
$ws.Cells.Item(12, 14).Value = 333.301639
$ws.Cells.Item(12, 15).Value = 0.3209442197221123
$ws.Cells.Item(12, 16).Value = 0.3209442197221123
$ws.Cells.Item(12, 17).Value = 10580.30019394972
$ws.Cells.Item(12, 18).Value = 95222.70174554751
$ws.Cells.Item(12, 19).Value = 0.06919448888761776
$ws.Cells.Item(12, 20).Value = 0.06919448888761776

# Row 13
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 95.23175666666667
$ws.Cells.Item(13, 8).Value = 285.69527
$ws.Cells.Item(13, 9).Value = 0.2155966197102082
$ws.Cells.Item(13, 10).Value = 0.2155966197102082
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 30.14303933333333
$ws.Cells.Item(13, 14).Value = 90.42911799999999
$ws.Cells.Item(13, 15).Value = 0.08707638763417187
$ws.Cells.Item(13, 16).Value = 0.08707638763417187
$ws.Cells.Item(13, 17).Value = 2870.574586985762
$ws.Cells.Item(13, 18).Value = 25835.17128287186
$ws.Cells.Item(13, 19).Value = 0.01877337483050323
$ws.Cells.Item(13, 20).Value = 0.01877337483050323

# Row 14
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 50.778675
$ws.Cells.Item(14, 8).Value = 152.336025
$ws.Cells.Item(14, 9).Value = 0.1149586132458188
$ws.Cells.Item(14, 10).Value = 0.1149586132458188
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 23.63579766666667
$ws.Cells.Item(14, 14).Value = 70.907393
$ws.Cells.Item(14, 15).Value = 0.06827844587621175
$ws.Cells.Item(14, 16).Value = 0.06827844587621175
$ws.Cells.Item(14, 17).Value = 1200.194488081425
$ws.Cells.Item(14, 18).Value = 10801.75039273283
$ws.Cells.Item(14, 19).Value = 0.007849195452508997
$ws.Cells.Item(14, 20).Value = 0.007849195452508997

# Row 15
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 50.778675
$ws.Cells.Item(15, 8).Value = 152.336025
$ws.Cells.Item(15, 9).Value = 0.1149586132458188
$ws.Cells.Item(15, 10).Value = 0.1149586132458188
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 181.2883913333334
$ws.Cells.Item(15, 14).Value = 543.865174
$ws.Cells.Item(15, 15).Value = 0.5237009467675041
$ws.Cells.Item(15, 16).Value = 0.523700946767504
$ws.Cells.Item(15, 17).Value = 9205.584304788152
$ws.Cells.Item(15, 18).Value = 82850.25874309336
$ws.Cells.Item(15, 19).Value = 0.06020393459591464
$ws.Cells.Item(15, 20).Value = 0.06020393459591463

# Row 16
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 50.778675
$ws.Cells.Item(16, 8).Value = 152.336025
$ws.Cells.Item(16, 9).Value = 0.1149586132458188
$ws.Cells.Item(16, 10).Value = 0.1149586132458188
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 111.1005463333333
$ws.Cells.Item(16, 14).Value = 333.301639
$ws.Cells.Item(16, 15).Value = 0.3209442197221123
$ws.Cells.Item(16, 16).Value = 0.3209442197221123
$ws.Cells.Item(16, 17).Value = 5641.538534582774
$ws.Cells.Item(16, 18).Value = 50773.84681124497
$ws.Cells.Item(16, 19).Value = 0.03689530242851539
$ws.Cells.Item(16, 20).Value = 0.03689530242851539

# Row 17
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 50.778675
$ws.Cells.Item(17, 8).Value = 152.336025
$ws.Cells.Item(17, 9).Value = 0.1149586132458188
$ws.Cells.Item(17, 10).Value = 0.1149586132458188
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 30.14303933333333
$ws.Cells.Item(17, 14).Value = 90.42911799999999
$ws.Cells.Item(17, 15).Value = 0.08707638763417187
$ws.Cells.Item(17, 16).Value = 0.08707638763417187
$ws.Cells.Item(17, 17).Value = 1530.62359781955
$ws.Cells.Item(17, 18).Value = 13775.61238037595
$ws.Cells.Item(17, 19).Value = 0.01001018076887976
$ws.Cells.Item(17, 20).Value = 0.01001018076887976
